$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data values between row 2 and row 4 for columns D, I, J, K, L, M, P
$ws.Range("D2").Value = 44370
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = 1080
$ws.Range("P2").Value = 180

$ws.Range("D4").Value = 44267
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 275
